$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H76").Value = 3146.1162
$ws.Range("I76").Value = 3105.3428
$ws.Range("J76").Value = 3324.5
$ws.Range("K76").Value = 3105.3428
$ws.Range("L76").Value = 3324.5
$ws.Range("M76").Value = -2790.3428
$ws.Range("N76").Value = -3954.5

$ws.Range("H79").Value = 3146.1162
$ws.Range("I79").Value = 3105.3428
$ws.Range("J79").Value = 3324.5
$ws.Range("K79").Value = 3105.3428
$ws.Range("L79").Value = 3324.5
$ws.Range("M79").Value = -2013.3428
$ws.Range("N79").Value = -5508.5

$ws.Range("H98").Value = 2099.8572
$ws.Range("I98").Value = 2116.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2116.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -618.5
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 2099.8572
$ws.Range("I122").Value = 2116.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6349.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3899.5
$ws.Range("N122").Value = -10900

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H31").Value = 5600
$ws.Range("I31").Value = 5600
$ws.Range("K31").Value = 5600
$ws.Range("M31").Value = -5306

$ws.Range("H32").Value = 18186494
$ws.Range("I32").Value = 22224994
$ws.Range("J32").Value = 13242.9
$ws.Range("K32").Value = 22224994
$ws.Range("L32").Value = 13242.9
$ws.Range("M32").Value = -22224707
$ws.Range("N32").Value = -13816.9

$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 10000
$ws.Range("K60").Value = 10000
$ws.Range("M60").Value = -9267

$ws.Range("H61").Value = 1685.3077
$ws.Range("I61").Value = 1389.6
$ws.Range("J61").Value = 2671
$ws.Range("K61").Value = 1389.6
$ws.Range("L61").Value = 2671
$ws.Range("M61").Value = -1177.6
$ws.Range("N61").Value = -3095

$ws.Range("H122").Value = 1246.12
$ws.Range("I122").Value = 1219.7084
$ws.Range("J122").Value = 1880
$ws.Range("K122").Value = 3659.1252
$ws.Range("L122").Value = 5640
$ws.Range("M122").Value = -1209.1252
$ws.Range("N122").Value = -10540

$ws.Range("H125").Value = 90048.336
$ws.Range("J125").Value = 90048.336
$ws.Range("L125").Value = 90048.336
$ws.Range("N125").Value = -99888.336

$ws.Range("H136").Value = 1685.3077
$ws.Range("I136").Value = 1389.6
$ws.Range("J136").Value = 2671
$ws.Range("K136").Value = 4168.799999999999
$ws.Range("L136").Value = 8013
$ws.Range("M136").Value = -1618.799999999999
$ws.Range("N136").Value = -13113

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0

$ws.Range("M24").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0

$ws.Range("N104").ClearContents()
$ws.Range("N109").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H21").Value = 550
$ws.Range("I21").Value = 550
$ws.Range("K21").Value = 1650
$ws.Range("M21").Value = -1477

$ws.Range("H131").Value = 706.8333
$ws.Range("I131").Value = 372.1111
$ws.Range("J131").Value = 907.6667
$ws.Range("K131").Value = 1116.3333
$ws.Range("L131").Value = 2723.0001
$ws.Range("M131").Value = 3923.6667
$ws.Range("N131").Value = -12803.0001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0

$ws.Range("H102").Value = 1620.4147
$ws.Range("I102").Value = 1622.8611
$ws.Range("J102").Value = 1602.8
$ws.Range("K102").Value = 1622.8611
$ws.Range("L102").Value = 1602.8
$ws.Range("M102").Value = -0.8611000000000786
$ws.Range("N102").Value = -4846.8

$ws.Range("H105").Value = 40756
$ws.Range("J105").Value = 40756
$ws.Range("L105").Value = 40756
$ws.Range("N105").Value = -47744

$ws.Range("H122").Value = 4763680.5
$ws.Range("I122").Value = 7693788.5
$ws.Range("J122").Value = 2254.875
$ws.Range("K122").Value = 23081365.5
$ws.Range("L122").Value = 6764.625
$ws.Range("M122").Value = -23078915.5
$ws.Range("N122").Value = -11664.625

$ws.Range("H126").Value = 4341.8945
$ws.Range("I126").Value = 4380
$ws.Range("J126").Value = 4314.1816
$ws.Range("K126").Value = 13140
$ws.Range("L126").Value = 12942.5448
$ws.Range("M126").Value = -10670
$ws.Range("N126").Value = -17882.5448

$ws.Range("H132").Value = 2606.04
$ws.Range("I132").Value = 2225.6667
$ws.Range("J132").Value = 4603
$ws.Range("K132").Value = 6677.000100000001
$ws.Range("L132").Value = 13809
$ws.Range("M132").Value = -4147.000100000001
$ws.Range("N132").Value = -18869

$ws.Range("N93").ClearContents()
$ws.Range("N98").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H106").Value = 23953.809
$ws.Range("J106").Value = 23953.809
$ws.Range("L106").Value = 23953.809
$ws.Range("N106").Value = -26477.809

$ws.Range("H122").Value = 3200
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -15100

$ws.Range("H132").Value = 3399.9395
$ws.Range("I132").Value = 3133.3333
$ws.Range("J132").Value = 3866.5
$ws.Range("K132").Value = 9399.999899999999
$ws.Range("L132").Value = 11599.5
$ws.Range("M132").Value = -6869.999899999999
$ws.Range("N132").Value = -16659.5005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H105").Value = 33333
$ws.Range("J105").Value = 33333
$ws.Range("L105").Value = 33333
$ws.Range("N105").Value = -40321

$ws.Range("H109").Value = 28530
$ws.Range("J109").Value = 28530
$ws.Range("L109").Value = 28530
$ws.Range("N109").Value = -31304

$ws.Range("H122").Value = 1765.6
$ws.Range("I122").Value = 1284.6818
$ws.Range("J122").Value = 2579.4614
$ws.Range("K122").Value = 3854.0454
$ws.Range("L122").Value = 7738.3842
$ws.Range("M122").Value = -1404.0454
$ws.Range("N122").Value = -12638.3842

$ws.Range("H126").Value = 2375.6316
$ws.Range("I126").Value = 2374.5557
$ws.Range("J126").Value = 2395
$ws.Range("K126").Value = 7123.6671
$ws.Range("L126").Value = 7185
$ws.Range("M126").Value = -4653.6671
$ws.Range("N126").Value = -12125

$ws.Range("H132").Value = 1412.3334
$ws.Range("I132").Value = 950.4194
$ws.Range("J132").Value = 3202.25
$ws.Range("K132").Value = 2851.2582
$ws.Range("L132").Value = 9606.75
$ws.Range("M132").Value = -321.2582000000002
$ws.Range("N132").Value = -14666.75
